# bug fix in Eduati data files
#
# Sheet1 ("SW620_noCTRL_meas" workbook) had 43 leftover junk rows
# (rows 45:87, column A only, sequential filler numbers) that don't belong
# with the real 14-column dataset (rows 1:44). Delete them, then restore
# the "as last saved" window/selection state: Sheet1 active (scrolled to
# around row 36, cell G58 selected) instead of Sheet3, and Sheet3 no
# longer marked as the selected tab.

$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the stray trailing rows (45:87), which also shrinks
#     the sheet's dimension from A1:N87 down to A1:N44 -------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
[void]$ws1.Activate()
[void]$ws1.Rows("45:87").Delete()

# Reproduce the saved view state for Sheet1: scrolled down so row 36 is
# the top visible row, with G58 as the active/selected cell.
[void]$ws1.Range("G58").Select()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet3 was the previously-active tab; simply re-activating Sheet1
#     (below) clears Sheet3's tabSelected flag on its own, leaving its
#     stored selection (A2:N44) untouched. ------------------------------

# Make Sheet1 the active sheet/tab again (also updates the workbook's
# bookViews activeTab so it points at Sheet1 instead of Sheet3).
[void]$ws1.Activate()
